$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("isa_template")

# Row 13: Tags
$ws1.Range("B13").Value = "Proteomics"
$ws1.Range("C13").Value = "Data Processing"
$ws1.Range("D13").Value = "Computation"
$ws1.Range("E13").Value = "software"
$ws1.Range("F13").Value = ""

# Row 14: Tags Term Accession Number
$ws1.Range("B14").Value = "http://purl.obolibrary.org/obo/NCIT_C20085"
$ws1.Range("C14").Value = "http://purl.obolibrary.org/obo/NCIT_C47925"
$ws1.Range("D14").Value = "http://purl.obolibrary.org/obo/NCIT_C61298"
$ws1.Range("E14").Value = "http://purl.obolibrary.org/obo/MS_1000531"
$ws1.Range("F14").Value = ""

# Row 15: Tags Term Source REF
$ws1.Range("B15").Value = "NCIT"
$ws1.Range("C15").Value = "NCIT"
$ws1.Range("D15").Value = "NCIT"
$ws1.Range("E15").Value = "MS"
$ws1.Range("F15").Value = ""
